$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.149.63'
$ws.Range("E2").Value = '  -2.02%  '

# Row 3
$ws.Range("D3").Value = '3.071.10'
$ws.Range("E3").Value = '  -2.76%  '

# Row 4
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.93%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.16%  '

# Row 7
$ws.Range("E7").Value = '  -0.09%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.552'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.88%  '

# Row 9
$ws.Range("D9").Value = '3.066.52'
$ws.Range("E9").Value = '  -2.66%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.155'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.88%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.85'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.37%  '

# Row 13
$ws.Range("E13").Value = '  -3.21%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.94%  '

# Row 15
$ws.Range("E15").Value = '  -1.92%  '

# Row 16
$ws.Range("D16").Value = '3.581.00'
$ws.Range("E16").Value = '  -2.82%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.01%  '

# Row 18
$ws.Range("D18").Value = '63.238.82'
$ws.Range("E18").Value = '  -1.50%  '

# Row 19
$ws.Range("D19").Value = '3.070.34'
$ws.Range("E19").Value = '  -2.82%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.09%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.87%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.716'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.77%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.51%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.99%  '

# Row 25
$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.36%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.06%  '

# Row 27
$ws.Range("E27").Value = '  -0.29%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.05%  '

# Row 29
$ws.Range("E29").Value = '  -0.04%  '

# Row 30
$ws.Range("E30").Value = '  -2.55%  '

# Row 31
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.41%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.19'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.86%  '

# Row 33
$ws.Range("E33").Value = '  +1.87%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.33%  '

# Row 35
$ws.Range("D35").Value = '0.0₃0847'
$ws.Range("E35").Value = '  +0.40%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.87%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.10'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.92%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.36%  '

# Row 39
$ws.Range("E39").Value = '  -4.91%  '

# Row 40
$ws.Range("E40").Value = '  +0.46%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.81%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '440.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.51%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.284'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.76%  '

# Row 44
$ws.Range("E44").Value = '  -3.86%  '

# Row 45
$ws.Range("D45").Value = '2.808.89'
$ws.Range("E45").Value = '  -4.38%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.52%  '

# Row 47
$ws.Range("E47").Value = '  +0.45%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.94%  '

# Row 50
$ws.Range("E50").Value = '  +0.71%  '

# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.80%  '
